$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 521.8182
$ws.Range("I2").Value = 177.14285
$ws.Range("J2").Value = 1125
$ws.Range("K2").Value = 177.14285
$ws.Range("L2").Value = 1125
$ws.Range("M2").Value = -64.14285000000001
$ws.Range("N2").Value = -1351
$ws.Range("H6").Value = 857.75
$ws.Range("I6").Value = 674.4
$ws.Range("J6").Value = 1163.3334
$ws.Range("K6").Value = 2023.2
$ws.Range("L6").Value = 3490.0002
$ws.Range("M6").Value = -1911.2
$ws.Range("N6").Value = -3714.0002
$ws.Range("H76").Value = 3009
$ws.Range("I76").Value = 2662.04
$ws.Range("K76").Value = 2662.04
$ws.Range("M76").Value = -2347.04
$ws.Range("H79").Value = 3009
$ws.Range("I79").Value = 2662.04
$ws.Range("K79").Value = 2662.04
$ws.Range("M79").Value = -1570.04
$ws.Range("H112").Value = 1147.0857
$ws.Range("J112").Value = 1119.0588
$ws.Range("L112").Value = 3357.1764
$ws.Range("N112").Value = -5573.1764
$ws.Range("H118").Value = 884
$ws.Range("I118").Value = 364
$ws.Range("K118").Value = 1092
$ws.Range("M118").Value = 565
$ws.Range("H127").Value = 812.8333
$ws.Range("I127").Value = 559.3570999999999
$ws.Range("J127").Value = 1700
$ws.Range("K127").Value = 1678.0713
$ws.Range("L127").Value = 5100
$ws.Range("M127").Value = 3281.9287
$ws.Range("N127").Value = -15020
$ws.Range("H135").Value = 406.975
$ws.Range("I135").Value = 396.08334
$ws.Range("J135").Value = 505
$ws.Range("K135").Value = 3564.75006
$ws.Range("L135").Value = 4545
$ws.Range("M135").Value = -1029.75006
$ws.Range("N135").Value = -9615
$ws.Range("H137").Value = 993.7708
$ws.Range("I137").Value = 911.5263
$ws.Range("J137").Value = 1306.3
$ws.Range("K137").Value = 2734.5789
$ws.Range("L137").Value = 3918.9
$ws.Range("M137").Value = -184.5789
$ws.Range("N137").Value = -9018.9
$ws.Range("H138").Value = 2918.0264
$ws.Range("I138").Value = 1495.0613
$ws.Range("J138").Value = 5500.4443
$ws.Range("K138").Value = 4485.1839
$ws.Range("L138").Value = 16501.3329
$ws.Range("M138").Value = 654.8161
$ws.Range("N138").Value = -26781.3329
$ws.Range("H141").Value = 4695.033
$ws.Range("I141").Value = 1149.9608
$ws.Range("J141").Value = 24783.777
$ws.Range("K141").Value = 3449.8824
$ws.Range("L141").Value = 74351.33099999999
$ws.Range("M141").Value = 1730.1176
$ws.Range("N141").Value = -84711.33099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6040.4546
$ws.Range("I28").Value = 3092.6
$ws.Range("J28").Value = 35519
$ws.Range("K28").Value = 3092.6
$ws.Range("L28").Value = 35519
$ws.Range("M28").Value = -2900.6
$ws.Range("N28").Value = -35903
$ws.Range("H32").Value = 3988.2727
$ws.Range("I32").Value = 2958.4106
$ws.Range("K32").Value = 2958.4106
$ws.Range("M32").Value = -2671.4106
$ws.Range("H41").Value = 36400
$ws.Range("I41").Value = 16200
$ws.Range("J41").Value = 46500
$ws.Range("K41").Value = 16200
$ws.Range("L41").Value = 46500
$ws.Range("M41").Value = -15786
$ws.Range("N41").Value = -47328
$ws.Range("H45").Value = 1393.2222
$ws.Range("I45").Value = 1091.6666
$ws.Range("J45").Value = 1996.3334
$ws.Range("K45").Value = 1091.6666
$ws.Range("L45").Value = 1996.3334
$ws.Range("M45").Value = -714.6666
$ws.Range("N45").Value = -2750.3334
$ws.Range("H61").Value = 1023.0469
$ws.Range("I61").Value = 581.13336
$ws.Range("K61").Value = 581.13336
$ws.Range("M61").Value = -369.13336
$ws.Range("H74").Value = 2906.5715
$ws.Range("I74").Value = 3162.6326
$ws.Range("J74").Value = 1114.1428
$ws.Range("K74").Value = 3162.6326
$ws.Range("L74").Value = 1114.1428
$ws.Range("M74").Value = -2288.6326
$ws.Range("N74").Value = -2862.1428
$ws.Range("H77").Value = 2906.5715
$ws.Range("I77").Value = 3162.6326
$ws.Range("J77").Value = 1114.1428
$ws.Range("K77").Value = 15813.163
$ws.Range("L77").Value = 5570.714
$ws.Range("M77").Value = -11445.163
$ws.Range("N77").Value = -14306.714
$ws.Range("H99").Value = 6040.4546
$ws.Range("I99").Value = 3092.6
$ws.Range("J99").Value = 35519
$ws.Range("K99").Value = 3092.6
$ws.Range("L99").Value = 35519
$ws.Range("M99").Value = -97.59999999999991
$ws.Range("N99").Value = -41509
$ws.Range("H122").Value = 1613.9354
$ws.Range("I122").Value = 1332.72
$ws.Range("J122").Value = 2785.6667
$ws.Range("K122").Value = 3998.16
$ws.Range("L122").Value = 8357.000100000001
$ws.Range("M122").Value = -1548.16
$ws.Range("N122").Value = -13257.0001
$ws.Range("H136").Value = 1023.0469
$ws.Range("I136").Value = 581.13336
$ws.Range("K136").Value = 1743.40008
$ws.Range("M136").Value = 806.5999199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 3000
$ws.Range("I11").Value = 3000
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -2860
$ws.Range("N11").Value = -3280
$ws.Range("H82").Value = 10611.111
$ws.Range("I82").Value = 1833.3334
$ws.Range("J82").Value = 15000
$ws.Range("K82").Value = 1833.3334
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = -1450.3334
$ws.Range("N82").Value = -15766
$ws.Range("H85").Value = 10611.111
$ws.Range("I85").Value = 1833.3334
$ws.Range("J85").Value = 15000
$ws.Range("K85").Value = 1833.3334
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = -507.3334
$ws.Range("N85").Value = -17652
$ws.Range("H99").Value = 890.5
$ws.Range("I99").Value = 788.125
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 788.125
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = 709.875
$ws.Range("N99").Value = -4296

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2335.2456
$ws.Range("I31").Value = 1640.5518
$ws.Range("J31").Value = 3054.75
$ws.Range("K31").Value = 1640.5518
$ws.Range("L31").Value = 3054.75
$ws.Range("M31").Value = -1345.5518
$ws.Range("N31").Value = -3644.75
$ws.Range("H34").Value = 2335.2456
$ws.Range("I34").Value = 1640.5518
$ws.Range("J34").Value = 3054.75
$ws.Range("K34").Value = 1640.5518
$ws.Range("L34").Value = 3054.75
$ws.Range("M34").Value = -1438.5518
$ws.Range("N34").Value = -3458.75
$ws.Range("H58").Value = 1305.3334
$ws.Range("I58").Value = 996.0172
$ws.Range("J58").Value = 2085.348
$ws.Range("K58").Value = 996.0172
$ws.Range("L58").Value = 2085.348
$ws.Range("M58").Value = -793.0172
$ws.Range("N58").Value = -2491.348
$ws.Range("H122").Value = 1984.3334
$ws.Range("I122").Value = 1230.2858
$ws.Range("J122").Value = 3040
$ws.Range("K122").Value = 3690.8574
$ws.Range("L122").Value = 9120
$ws.Range("M122").Value = -1240.8574
$ws.Range("N122").Value = -14020
$ws.Range("H132").Value = 1275.2778
$ws.Range("I132").Value = 604.0244
$ws.Range("J132").Value = 3392.3076
$ws.Range("K132").Value = 1812.0732
$ws.Range("L132").Value = 10176.9228
$ws.Range("M132").Value = 717.9268
$ws.Range("N132").Value = -15236.9228
$ws.Range("H134").Value = 1441.3226
$ws.Range("I134").Value = 1331.38
$ws.Range("J134").Value = 1899.4166
$ws.Range("K134").Value = 3994.14
$ws.Range("L134").Value = 5698.2498
$ws.Range("M134").Value = -1459.14
$ws.Range("N134").Value = -10768.2498
$ws.Range("H136").Value = 1305.3334
$ws.Range("I136").Value = 996.0172
$ws.Range("J136").Value = 2085.348
$ws.Range("K136").Value = 2988.0516
$ws.Range("L136").Value = 6256.044
$ws.Range("M136").Value = -438.0515999999998
$ws.Range("N136").Value = -11356.044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2134.6785
$ws.Range("I122").Value = 1792.5294
$ws.Range("K122").Value = 5377.5882
$ws.Range("M122").Value = -2927.5882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10105037
$ws.Range("I122").Value = 15875118
$ws.Range("J122").Value = 7395
$ws.Range("K122").Value = 47625354
$ws.Range("L122").Value = 22185
$ws.Range("M122").Value = -47622904
$ws.Range("N122").Value = -27085
$ws.Range("H132").Value = 4387.6665
$ws.Range("I132").Value = 3983.6711
$ws.Range("J132").Value = 5862.25
$ws.Range("K132").Value = 11951.0133
$ws.Range("L132").Value = 17586.75
$ws.Range("M132").Value = -9421.013300000001
$ws.Range("N132").Value = -22646.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 73853.78999999999
$ws.Range("I122").Value = 127307.25
$ws.Range("J122").Value = 2582.5
$ws.Range("K122").Value = 381921.75
$ws.Range("L122").Value = 7747.5
$ws.Range("M122").Value = -379471.75
$ws.Range("N122").Value = -12647.5
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120
